# Fix species names: trim trailing whitespace from "aleppo pine " and "chinese elm "
# Rows 126-187 contain "aleppo pine " -> "aleppo pine"
# Rows 250-311 contain "chinese elm " -> "chinese elm"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A250:A311").Value = "chinese elm"
$ws.Range("A126:A187").Value = "aleppo pine"

# Update sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 237
$ws.Range("A251:A311").Select()
